$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1567.0819
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1567.0819
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4701.245699999999
$ws.Range("N17").Value = -5037.245699999999
$ws.Range("H64").Value = 3054
$ws.Range("I64").Value = 3096
$ws.Range("J64").Value = 2949
$ws.Range("K64").Value = 3096
$ws.Range("L64").Value = 2949
$ws.Range("M64").Value = -2848
$ws.Range("N64").Value = -3445
$ws.Range("H67").Value = 3054
$ws.Range("I67").Value = 3096
$ws.Range("J67").Value = 2949
$ws.Range("K67").Value = 3096
$ws.Range("L67").Value = 2949
$ws.Range("M67").Value = -2238
$ws.Range("N67").Value = -4665
$ws.Range("H88").Value = 5334.407
$ws.Range("I88").Value = 4497.2856
$ws.Range("J88").Value = 6235.923
$ws.Range("K88").Value = 4497.2856
$ws.Range("L88").Value = 6235.923
$ws.Range("M88").Value = -4091.2856
$ws.Range("N88").Value = -7047.923
$ws.Range("H91").Value = 5334.407
$ws.Range("I91").Value = 4497.2856
$ws.Range("J91").Value = 6235.923
$ws.Range("K91").Value = 4497.2856
$ws.Range("L91").Value = 6235.923
$ws.Range("M91").Value = -3093.2856
$ws.Range("N91").Value = -9043.922999999999
$ws.Range("H115").Value = 1477.5454
$ws.Range("I115").Value = 650.6
$ws.Range("J115").Value = 2166.6667
$ws.Range("K115").Value = 1951.8
$ws.Range("L115").Value = 6500.000100000001
$ws.Range("M115").Value = -384.8000000000002
$ws.Range("N115").Value = -9634.000100000001
$ws.Range("H125").Value = 2461.818
$ws.Range("I125").Value = 2325
$ws.Range("J125").Value = 2540
$ws.Range("K125").Value = 20925
$ws.Range("L125").Value = 22860
$ws.Range("M125").Value = -18465
$ws.Range("N125").Value = -27780
$ws.Range("H131").Value = 4167.067
$ws.Range("I131").Value = 3969.375
$ws.Range("J131").Value = 4393
$ws.Range("K131").Value = 11908.125
$ws.Range("L131").Value = 13179
$ws.Range("M131").Value = -6868.125
$ws.Range("N131").Value = -23259
$ws.Range("H135").Value = 761.8333
$ws.Range("I135").Value = 535.7381
$ws.Range("J135").Value = 1553.1666
$ws.Range("K135").Value = 4821.642900000001
$ws.Range("L135").Value = 13978.4994
$ws.Range("M135").Value = -2286.642900000001
$ws.Range("N135").Value = -19048.4994
$ws.Range("H141").Value = 320956.38
$ws.Range("I141").Value = 1321.1072
$ws.Range("J141").Value = 1439679.9
$ws.Range("K141").Value = 3963.3216
$ws.Range("L141").Value = 4319039.699999999
$ws.Range("M141").Value = 1216.6784
$ws.Range("N141").Value = -4329399.699999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1993.4642
$ws.Range("I122").Value = 1585.9048
$ws.Range("J122").Value = 3216.1428
$ws.Range("K122").Value = 4757.7144
$ws.Range("L122").Value = 9648.428400000001
$ws.Range("M122").Value = -2307.7144
$ws.Range("N122").Value = -14548.4284
$ws.Range("H132").Value = 2133.7551
$ws.Range("I132").Value = 1612.4286
$ws.Range("J132").Value = 2828.8572
$ws.Range("K132").Value = 4837.2858
$ws.Range("L132").Value = 8486.571599999999
$ws.Range("M132").Value = -2307.2858
$ws.Range("N132").Value = -13546.5716
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5134.9
$ws.Range("I20").Value = 5915.385
$ws.Range("J20").Value = 3685.4285
$ws.Range("K20").Value = 5915.385
$ws.Range("L20").Value = 3685.4285
$ws.Range("M20").Value = -5668.385
$ws.Range("N20").Value = -4179.4285
$ws.Range("H132").Value = 29500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 29500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 29500
$ws.Range("H134").Value = 2322.9736
$ws.Range("I134").Value = 1837.0344
$ws.Range("J134").Value = 3888.7778
$ws.Range("K134").Value = 5511.1032
$ws.Range("L134").Value = 11666.3334
$ws.Range("M134").Value = -2976.1032
$ws.Range("N134").Value = -16736.3334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1779
$ws.Range("I22").Value = 603.1667
$ws.Range("J22").Value = 3190
$ws.Range("K22").Value = 603.1667
$ws.Range("L22").Value = 3190
$ws.Range("M22").Value = -253.1667
$ws.Range("N22").Value = -3890
$ws.Range("H62").Value = 3089.4443
$ws.Range("I62").Value = 2922
$ws.Range("J62").Value = 3153.8462
$ws.Range("K62").Value = 2922
$ws.Range("L62").Value = 3153.8462
$ws.Range("M62").Value = -2298
$ws.Range("N62").Value = -4401.8462
$ws.Range("H65").Value = 3089.4443
$ws.Range("I65").Value = 2922
$ws.Range("J65").Value = 3153.8462
$ws.Range("K65").Value = 14610
$ws.Range("L65").Value = 15769.231
$ws.Range("M65").Value = -11490
$ws.Range("N65").Value = -22009.231
$ws.Range("H99").Value = 2599.8
$ws.Range("I99").Value = 1516.6666
$ws.Range("J99").Value = 4224.5
$ws.Range("K99").Value = 1516.6666
$ws.Range("L99").Value = 4224.5
$ws.Range("M99").Value = -18.66660000000002
$ws.Range("N99").Value = -7220.5
$ws.Range("H126").Value = 2599.8
$ws.Range("I126").Value = 1516.6666
$ws.Range("J126").Value = 4224.5
$ws.Range("K126").Value = 4549.9998
$ws.Range("L126").Value = 12673.5
$ws.Range("M126").Value = -2079.9998
$ws.Range("N126").Value = -17613.5
$ws.Range("H127").Value = 33000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 33000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 33000
$ws.Range("N127").Value = -42920
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2790.56
$ws.Range("I3").Value = 1678.6
$ws.Range("J3").Value = 4458.5
$ws.Range("K3").Value = 5035.799999999999
$ws.Range("L3").Value = 13375.5
$ws.Range("M3").Value = -4923.799999999999
$ws.Range("N3").Value = -13599.5
$ws.Range("H59").Value = 2440
$ws.Range("I59").Value = 300
$ws.Range("J59").Value = 2975
$ws.Range("K59").Value = 900
$ws.Range("L59").Value = 8925
$ws.Range("N59").Value = -10005
$ws.Range("H64").Value = 2100.9285
$ws.Range("I64").Value = 1244.7142
$ws.Range("J64").Value = 2957.1428
$ws.Range("K64").Value = 3734.1426
$ws.Range("L64").Value = 8871.428400000001
$ws.Range("M64").Value = -3464.1426
$ws.Range("N64").Value = -9411.428400000001
$ws.Range("H67").Value = 2100.9285
$ws.Range("I67").Value = 1244.7142
$ws.Range("J67").Value = 2957.1428
$ws.Range("K67").Value = 3734.1426
$ws.Range("L67").Value = 8871.428400000001
$ws.Range("M67").Value = -2798.1426
$ws.Range("N67").Value = -10743.4284
$ws.Range("H99").Value = 1751.9231
$ws.Range("I99").Value = 1462.5
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 4387.5
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -2141.5
$ws.Range("N99").Value = -10492
$ws.Range("H116").Value = 1902.75
$ws.Range("I116").Value = 1203.6666
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 3610.9998
$ws.Range("L116").Value = 12000
$ws.Range("M116").Value = -168.9998000000001
$ws.Range("N116").Value = -18884
$ws.Range("H120").Value = 17009.46
$ws.Range("I120").Value = 13522.5
$ws.Range("J120").Value = 18559.223
$ws.Range("K120").Value = 40567.5
$ws.Range("L120").Value = 55677.66900000001
$ws.Range("M120").Value = -35729.5
$ws.Range("N120").Value = -65353.66900000001
$ws.Range("H122").Value = 1311.091
$ws.Range("I122").Value = 385
$ws.Range("J122").Value = 1840.2858
$ws.Range("K122").Value = 3465
$ws.Range("L122").Value = 16562.5722
$ws.Range("M122").Value = -1015
$ws.Range("N122").Value = -21462.5722
$ws.Range("H138").Value = 4669.1
$ws.Range("I138").Value = 1485.25
$ws.Range("J138").Value = 6791.6665
$ws.Range("K138").Value = 4455.75
$ws.Range("L138").Value = 20374.9995
$ws.Range("M138").Value = 684.25
$ws.Range("N138").Value = -30654.9995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1618.8448
$ws.Range("I132").Value = 973.75
$ws.Range("J132").Value = 2674.4546
$ws.Range("K132").Value = 2921.25
$ws.Range("L132").Value = 8023.3638
$ws.Range("M132").Value = -391.25
$ws.Range("N132").Value = -13083.3638
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 71431896
$ws.Range("I61").Value = 111113950
$ws.Range("J61").Value = 4200
$ws.Range("K61").Value = 111113950
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -111113748
$ws.Range("N61").Value = -4604
$ws.Range("H113").Value = 71431896
$ws.Range("I113").Value = 111113950
$ws.Range("J113").Value = 4200
$ws.Range("K113").Value = 111113950
$ws.Range("L113").Value = 4200
$ws.Range("M113").Value = -111111780
$ws.Range("N113").Value = -8540
$ws.Range("H122").Value = 3061.4644
$ws.Range("I122").Value = 2642.8096
$ws.Range("J122").Value = 4317.4287
$ws.Range("K122").Value = 7928.4288
$ws.Range("L122").Value = 12952.2861
$ws.Range("M122").Value = -5478.4288
$ws.Range("N122").Value = -17852.2861
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2511395
$ws.Range("I3").Value = 3333526.8
$ws.Range("J3").Value = 45000
$ws.Range("K3").Value = 3333526.8
$ws.Range("L3").Value = 45000
$ws.Range("M3").Value = -3333412.8
$ws.Range("N3").Value = -45228
$ws.Range("H46").Value = 48809.668
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 48809.668
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 48809.668
$ws.Range("N46").Value = -49271.668
$ws.Range("H122").Value = 386748.06
$ws.Range("I122").Value = 477948.16
$ws.Range("J122").Value = 3707.8
$ws.Range("K122").Value = 1433844.48
$ws.Range("L122").Value = 11123.4
$ws.Range("M122").Value = -1431394.48
$ws.Range("N122").Value = -16023.4
$ws.Range("H132").Value = 16087.513
$ws.Range("I132").Value = 3385.8928
$ws.Range("J132").Value = 48418.91
$ws.Range("K132").Value = 10157.6784
$ws.Range("L132").Value = 145256.73
$ws.Range("M132").Value = -7627.678400000001
$ws.Range("N132").Value = -150316.73
$ws.Range("H134").Value = 48809.668
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 48809.668
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 146429.004
$ws.Range("N134").Value = -151499.004
